$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Data" sheet: append 8 new daily observations (rows 451-458), matching
#    the formatting already used by the preceding row (450).
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Clone the formatting of the last existing data row onto the new rows so the
# new date cells keep the same date-number-format / border / font style as
# every other row in column A (and column B stays the default style).
$wsData.Range("A450:B450").Copy($wsData.Range("A451:B458"))

$newDates = @(45131, 45132, 45133, 45134, 45135, 45138, 45139, 45140)
$newValues = @(1770.867, 1720.716, 1749.733, 1735.783, 1730.227, 1821.124, 1739.554, 1770.186)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = 451 + $i
    $wsData.Cells.Item($row, 1).Value = $newDates[$i]
    $wsData.Cells.Item($row, 2).Value = $newValues[$i]
}

# ---------------------------------------------------------------------------
# 2) "SeriesInfo" sheet: refresh the metadata that FRED stamps on every pull
#    (realtime_start / realtime_end / observation_end / last_updated).
#    These must stay plain text cells (as they were before). The new text
#    ("2023-08-03", "2023-08-02", ...) looks exactly like a date, so Excel's
#    auto-detection would otherwise silently convert the cell into a
#    date-typed number. Force the cell to Text format *before* writing the
#    value so it is stored verbatim, then restore the cell's original
#    (default/no explicit style) appearance by copying the style from a
#    neighboring plain-text cell.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")
$plainStyle = $wsInfo.Range("A2").Style

function Set-PlainTextValue($range, $value, $styleSource) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $styleSource
}

Set-PlainTextValue $wsInfo.Range("B3") "2023-08-03" $plainStyle
Set-PlainTextValue $wsInfo.Range("B4") "2023-08-03" $plainStyle
Set-PlainTextValue $wsInfo.Range("B7") "2023-08-02" $plainStyle
Set-PlainTextValue $wsInfo.Range("B14") "2023-08-02 13:01:05-05" $plainStyle
